$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Microgrid" survey rows (A20:B25) had only Source/Survey labels with no
# household counts yet -- fill in the missing "Households" column values.
$ws.Range("C20").Value = 20
$ws.Range("C21").Value = 10
$ws.Range("C22").Value = 0
$ws.Range("C23").Value = 0
$ws.Range("C24").Value = 13
$ws.Range("C25").Value = 75

# Reflect where the editing session left the cursor/selection.
$ws.Range("E26").Select()
